$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coinranking "Price" values are scraped as plain text (note the
# thousand-separator dots, e.g. "68.798.81"). Several of the updated
# prices happen to look like ordinary decimals (e.g. "1.00", "601.77"),
# and a bare Range.Value assignment would let Excel reinterpret those
# as numbers. Briefly format the cell as Text ("@") while writing the
# value, then restore the default "Normal" style so the cell keeps
# its original (unstyled) look.

$ws.Range("D2").Value = "68.798.81"
$ws.Range("E2").Value = "  +2.42%  "
$ws.Range("D3").Value = "3.748.02"
$ws.Range("E3").Value = "  +1.95%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "601.77"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.45%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "168.73"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.73%  "
$ws.Range("D7").Value = "3.746.53"
$ws.Range("E7").Value = "  +1.93%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").Value = "  +2.33%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.165"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.44%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.35"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.30%  "
$ws.Range("E12").Value = "  +0.48%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "38.33"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.98%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000248"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.70%  "
$ws.Range("D15").Value = "4.374.78"
$ws.Range("E15").Value = "  +2.01%  "
$ws.Range("D16").Value = "3.748.43"
$ws.Range("E16").Value = "  +2.04%  "
$ws.Range("D17").Value = "68.800.20"
$ws.Range("E17").Value = "  +2.33%  "
$ws.Range("E18").Value = "  +2.64%  "
$ws.Range("E19").Value = "  +1.26%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.08"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.40%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.81"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +19.04%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "496.27"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.62%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.728"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.47%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000154"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +10.04%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "85.46"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.26%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.32"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.10%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.32"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.52%  "
$ws.Range("E28").Value = "  +3.73%  "
$ws.Range("E29").Value = "  +0.43%  "
$ws.Range("E30").Value = "  +2.69%  "
$ws.Range("E31").Value = "  +6.77%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.91"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.60%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "31.77"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.28%  "
$ws.Range("D34").Value = "3.892.73"
$ws.Range("E34").Value = "  +2.13%  "
$ws.Range("D35").Value = "3.682.96"
$ws.Range("E35").Value = "  +1.94%  "
$ws.Range("E36").Value = "  +1.93%  "
$ws.Range("E37").Value = "  +0.02%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.02"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.85%  "
$ws.Range("E39").Value = "  +1.93%  "
$ws.Range("E40").Value = "  +1.60%  "
$ws.Range("E41").Value = "  +0.65%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "438.78"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.29%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.94"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +6.17%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "48.83"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.62%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.97"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.26%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.47"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.91%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "40.43"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.03%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "141.82"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.07%  "
$ws.Range("D50").Value = "2.788.92"
$ws.Range("E50").Value = "  +1.40%  "
$ws.Range("E51").Value = "  +2.78%  "
